# StomatalDensity.xlsx - "moved some files/worked on stomata"
#
# Fill in measured stomatal-density data (B = area, C = count, D = C/B
# computed by the existing shared formula) for the PIPO13/14/15/16 tree
# groups on the PIPO sheet, expanding the PIPO14/15/16 groups from a
# single placeholder row each to three measured rows each (inserting two
# new rows per group), and updating the view (zoom / frozen-pane scroll /
# selection) to where the user ended up after the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIPO")

# --- PIPO13 group (rows 38-40 already exist, just needs values) ---
$ws.Range("B38").Value = 1.732
$ws.Range("C38").Value = 91

$ws.Range("B39").Value = 1.1200000000000001
$ws.Range("C39").Value = 47

$ws.Range("B40").Value = 0.35599999999999998
$ws.Range("C40").Value = 19

# --- PIPO14 group: first row (41) already exists; add two more rows ---
$ws.Range("B41").Value = 1.08
$ws.Range("C41").Value = 83

$ws.Rows("42:43").Insert()

$ws.Range("A42").Value = "PIPO14"
$ws.Range("B42").Value = 0.72699999999999998
$ws.Range("C42").Value = 59
$ws.Range("D42").Formula = "=C42/B42"

$ws.Range("A43").Value = "PIPO14"
$ws.Range("B43").Value = 0.24099999999999999
$ws.Range("C43").Value = 13
$ws.Range("D43").Formula = "=C43/B43"

# --- PIPO15 group: first row (now 44, was 42) already exists; add two more ---
$ws.Range("B44").Value = 0.66600000000000004
$ws.Range("C44").Value = 49

$ws.Rows("45:46").Insert()

$ws.Range("A45").Value = "PIPO15"
$ws.Range("B45").Value = 0.49
$ws.Range("C45").Value = 42
$ws.Range("D45").Formula = "=C45/B45"

$ws.Range("A46").Value = "PIPO15"
$ws.Range("B46").Value = 1.149
$ws.Range("C46").Value = 82
$ws.Range("D46").Formula = "=C46/B46"

# --- PIPO16 group: first row (now 47, was 43) already exists; add two more ---
$ws.Range("B47").Value = 2.468
$ws.Range("C47").Value = 134

$ws.Rows("48:49").Insert()

$ws.Range("A48").Value = "PIPO16"
$ws.Range("B48").Value = 1.264
$ws.Range("C48").Value = 75
$ws.Range("D48").Formula = "=C48/B48"

$ws.Range("A49").Value = "PIPO16"
$ws.Range("B49").Value = 0.82799999999999996
$ws.Range("C49").Value = 50
$ws.Range("D49").Formula = "=C49/B49"

# --- Highlight (bold) every row that now has real data, rows 38-49 ---
$ws.Range("A38:A49").Font.Bold = $true

# --- View state: zoom, frozen-pane scroll position, selection ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 188

$ws.Range("A47:A49").Select()
$ws.Range("A49").Activate()

$win.ScrollRow = 39
$win.ScrollColumn = 1
